$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (Mercredi, week "Semaine 18 au 24 Mai") ---
# C16: Pause range (time-range text, same formatting family as other "Pause" cells)
$ws.Range("B16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = "12:47 - 15:00"

# D16: Fin time (time number format, like the rest of column D)
$ws.Range("B16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 0.82638888888888884

# E16: Temps total (plain centered text format)
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = "9h09"

# --- Row 17 (Jeudi) ---
# B17: Début time
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = 0.36874999999999997

$excel.CutCopyMode = 0

# Update active selection to C17
$ws.Range("C17").Select()
